# Update "想去人数" (interested-count) values in column F on two sheets,
# matching the data refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st tab)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1178
$ws1.Range("F4").Value = 0
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 9273
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 522
$ws1.Range("F12").Value = 74

# Sheet "全部类型" (4th tab)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14
$ws4.Range("F3").Value = 640
$ws4.Range("F8").Value = 522
$ws4.Range("F10").Value = 9273
$ws4.Range("F11").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 0
